# Auto update Excel log
# Appends new sensor-reading rows captured on 2026-01-28 to the PIR,
# Humidity and Temperature sheets of the master log workbook.

$wb = $excel.ActiveWorkbook

# --- PIR sheet: new rows 173-186 ---
$ws = $wb.Worksheets.Item("PIR")

$PIRRows = @(
    ,@('2026-01-28', '16:24:54', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:24:56', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:24:57', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:25:03', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:25:07', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:25:12', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:25:17', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:25:23', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:25:28', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:25:33', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:25:38', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:25:43', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:25:48', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:25:53', '16:00', 'Bathroom', 'No Motion', 'Inactive')
)

# Column A holds a literal date-like string (e.g. 2026-01-28); pre-format
# it as Text so Excel does not silently convert it to a date serial value.
$ws.Range("A173:A186").NumberFormat = "@"

$r = 173
foreach ($row in $PIRRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Humidity sheet: new rows 173-184 ---
$ws = $wb.Worksheets.Item("Humidity")

$HumidityRows = @(
    ,@('2026-01-28', '16:24:55', '16:00', 'Bathroom', '88.3%', 'Active')
    ,@('2026-01-28', '16:24:56', '16:00', 'Bathroom', '87.4%', 'Active')
    ,@('2026-01-28', '16:25:02', '16:00', 'Bathroom', '88.3%', 'Active')
    ,@('2026-01-28', '16:25:06', '16:00', 'Bathroom', '87.4%', 'Active')
    ,@('2026-01-28', '16:25:10', '16:00', 'Bathroom', '88.3%', 'Active')
    ,@('2026-01-28', '16:25:14', '16:00', 'Bathroom', '87.4%', 'Active')
    ,@('2026-01-28', '16:25:18', '16:00', 'Bathroom', '88.3%', 'Active')
    ,@('2026-01-28', '16:25:22', '16:00', 'Bathroom', '88.3%', 'Active')
    ,@('2026-01-28', '16:25:30', '16:00', 'Bathroom', '88.3%', 'Active')
    ,@('2026-01-28', '16:25:42', '16:00', 'Bathroom', '88.3%', 'Active')
    ,@('2026-01-28', '16:25:46', '16:00', 'Bathroom', '87.4%', 'Active')
    ,@('2026-01-28', '16:25:50', '16:00', 'Bathroom', '88.3%', 'Active')
)

# Column A holds a literal date-like string (e.g. 2026-01-28); pre-format
# it as Text so Excel does not silently convert it to a date serial value.
$ws.Range("A173:A184").NumberFormat = "@"
# Column E holds a literal percentage-like string (e.g. 88.3%); pre-format
# it as Text so Excel does not convert it to a numeric percentage value.
$ws.Range("E173:E184").NumberFormat = "@"

$r = 173
foreach ($row in $HumidityRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Temperature sheet: new rows 173-184 ---
$ws = $wb.Worksheets.Item("Temperature")

$TemperatureRows = @(
    ,@('2026-01-28', '16:24:55', '16:00', 'Bathroom', '22.7C', 'Active')
    ,@('2026-01-28', '16:24:57', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:25:02', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:25:06', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:25:10', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:25:14', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:25:18', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:25:22', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:25:30', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:25:42', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:25:46', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:25:51', '16:00', 'Bathroom', '22.8C', 'Active')
)

# Column A holds a literal date-like string (e.g. 2026-01-28); pre-format
# it as Text so Excel does not silently convert it to a date serial value.
$ws.Range("A173:A184").NumberFormat = "@"

$r = 173
foreach ($row in $TemperatureRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

